$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = "sd"
$ws.Range("J8").Value = "Statement-non-opinion"
$ws.Range("I26").Value = "aa"
$ws.Range("J26").Value = "Agree/Accept"
$ws.Range("I27").Value = "sd"
$ws.Range("J27").Value = "Statement-non-opinion"
$ws.Range("I43").Value = "ba"
$ws.Range("J43").Value = "Appreciation"
$ws.Range("I50").Value = "sd"
$ws.Range("J50").Value = "Statement-non-opinion"
$ws.Range("I55").Value = "sd"
$ws.Range("J55").Value = "Statement-non-opinion"
$ws.Range("I57").Value = "aa"
$ws.Range("J57").Value = "Agree/Accept"
$ws.Range("I64").Value = "ba"
$ws.Range("J64").Value = "Appreciation"
$ws.Range("I66").Value = "ba"
$ws.Range("J66").Value = "Appreciation"
$ws.Range("I80").Value = "sd"
$ws.Range("J80").Value = "Statement-non-opinion"
$ws.Range("I85").Value = "sd"
$ws.Range("J85").Value = "Statement-non-opinion"
$ws.Range("I86").Value = "sd"
$ws.Range("J86").Value = "Statement-non-opinion"
$ws.Range("I92").Value = "sv"
$ws.Range("J92").Value = "Statement-opinion"
$ws.Range("I96").Value = "aa"
$ws.Range("J96").Value = "Agree/Accept"
$ws.Range("I98").Value = "ba"
$ws.Range("J98").Value = "Appreciation"
$ws.Range("I100").Value = "b"
$ws.Range("J100").Value = "Acknowledge (Backchannel)"
$ws.Range("I111").Value = "b"
$ws.Range("J111").Value = "Acknowledge (Backchannel)"
$ws.Range("I116").Value = "aa"
$ws.Range("J116").Value = "Agree/Accept"
$ws.Range("I118").Value = "b"
$ws.Range("J118").Value = "Acknowledge (Backchannel)"
$ws.Range("I120").Value = "sd"
$ws.Range("J120").Value = "Statement-non-opinion"
$ws.Range("I121").Value = "sd"
$ws.Range("J121").Value = "Statement-non-opinion"
$ws.Range("I122").Value = "%"
$ws.Range("J122").Value = "Uninterpretable"
$ws.Range("I125").Value = "sd"
$ws.Range("J125").Value = "Statement-non-opinion"
$ws.Range("I138").Value = "sd"
$ws.Range("J138").Value = "Statement-non-opinion"
$ws.Range("I144").Value = "sd"
$ws.Range("J144").Value = "Statement-non-opinion"
$ws.Range("I148").Value = "b"
$ws.Range("J148").Value = "Acknowledge (Backchannel)"
$ws.Range("I149").Value = "ba"
$ws.Range("J149").Value = "Appreciation"
$ws.Range("I152").Value = "%"
$ws.Range("J152").Value = "Uninterpretable"
$ws.Range("I157").Value = "aa"
$ws.Range("J157").Value = "Agree/Accept"
$ws.Range("I158").Value = "aa"
$ws.Range("J158").Value = "Agree/Accept"
$ws.Range("I160").Value = "b"
$ws.Range("J160").Value = "Acknowledge (Backchannel)"
$ws.Range("I162").Value = "aa"
$ws.Range("J162").Value = "Agree/Accept"
$ws.Range("I196").Value = "b"
$ws.Range("J196").Value = "Acknowledge (Backchannel)"
$ws.Range("I205").Value = "sv"
$ws.Range("J205").Value = "Statement-opinion"
$ws.Range("I208").Value = "b"
$ws.Range("J208").Value = "Acknowledge (Backchannel)"
$ws.Range("I226").Value = "b"
$ws.Range("J226").Value = "Acknowledge (Backchannel)"
$ws.Range("I234").Value = "b"
$ws.Range("J234").Value = "Acknowledge (Backchannel)"
